$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 3.6
$ws.Range("K2").Value = 9.5
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 1.73
$ws.Range("T2").Value = 6
$ws.Range("U2").Value = 7.5
$ws.Range("Z2").Value = 9.5
$ws.Range("AA2").Value = 7.5
$ws.Range("AB2").Value = 19
$ws.Range("AC2").Value = 67
$ws.Range("AE2").Value = 11
$ws.Range("G5").Value = 1.6
$ws.Range("H5").Value = 4.2
$ws.Range("I5").Value = 5.25
$ws.Range("J5").Value = 1.04
$ws.Range("K5").Value = 13
$ws.Range("N5").Value = 1.8
$ws.Range("O5").Value = 2
$ws.Range("AA5").Value = 8
$ws.Range("AD5").Value = 301
$ws.Range("G6").Value = 1.36
$ws.Range("H6").Value = 4.65
$ws.Range("I6").Value = 6.3
$ws.Range("T6").Value = 9.25
$ws.Range("U6").Value = 7.6
$ws.Range("V6").Value = 7.3
$ws.Range("X6").Value = 8.5
$ws.Range("Y6").Value = 15
$ws.Range("Z6").Value = 20
$ws.Range("AA6").Value = 9
$ws.Range("AE6").Value = 22
$ws.Range("AH6").Value = 100
$ws.Range("AJ6").Value = 32
$ws.Range("G7").Value = 1.25
$ws.Range("H7").Value = 4.8
$ws.Range("I7").Value = 10
$ws.Range("N7").Value = 1.49
$ws.Range("O7").Value = 2.44
$ws.Range("R7").Value = 1.86
$ws.Range("S7").Value = 1.85
$ws.Range("T7").Value = 6.8
$ws.Range("U7").Value = 5.7
$ws.Range("V7").Value = 7.3
$ws.Range("W7").Value = 6.7
$ws.Range("X7").Value = 8.75
$ws.Range("Y7").Value = 19.5
$ws.Range("Z7").Value = 14.5
$ws.Range("AA7").Value = 8.75
$ws.Range("AB7").Value = 16.5
$ws.Range("AC7").Value = 60
$ws.Range("AD7").Value = 350
$ws.Range("AE7").Value = 25
$ws.Range("AF7").Value = 65
$ws.Range("AG7").Value = 25
$ws.Range("AH7").Value = 250
$ws.Range("AI7").Value = 90
$ws.Range("AJ7").Value = 60
$ws.Range("G8").Value = 1.91
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 3.5
$ws.Range("N8").Value = 1.62
$ws.Range("O8").Value = 2.25
$ws.Range("U8").Value = 11
$ws.Range("AC8").Value = 41
$ws.Range("AG8").Value = 13
$ws.Range("AI8").Value = 26
$ws.Range("G9").Value = 1.75
$ws.Range("H9").Value = 3.9
$ws.Range("I9").Value = 3.8
$ws.Range("L9").Value = 1.25
$ws.Range("M9").Value = 3.7
$ws.Range("P9").Value = 1.36
$ws.Range("Q9").Value = 2.95
$ws.Range("X9").Value = 13
$ws.Range("AC9").Value = 51
$ws.Range("AI9").Value = 34
$ws.Range("H10").Value = 3.9
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 1.03
$ws.Range("K10").Value = 9.75
$ws.Range("L10").Value = 1.18
$ws.Range("M10").Value = 4.45
$ws.Range("N10").Value = 1.57
$ws.Range("O10").Value = 2.35
$ws.Range("P10").Value = 1.3
$ws.Range("Q10").Value = 3.3
$ws.Range("R10").Value = 1.57
$ws.Range("S10").Value = 2.25
$ws.Range("T10").Value = 10
$ws.Range("U10").Value = 10
$ws.Range("W10").Value = 15
$ws.Range("Z10").Value = 17
$ws.Range("AB10").Value = 13
$ws.Range("AD10").Value = 126
$ws.Range("AE10").Value = 17
$ws.Range("AG10").Value = 13
$ws.Range("AI10").Value = 29
$ws.Range("AJ10").Value = 29
$ws.Range("G12").Value = 2.2
$ws.Range("I12").Value = 3.1
$ws.Range("N12").Value = 2
$ws.Range("O12").Value = 1.85
$ws.Range("P12").Value = 1.44
$ws.Range("Q12").Value = 2.67
$ws.Range("W12").Value = 21
$ws.Range("Z12").Value = 10
$ws.Range("AB12").Value = 13
$ws.Range("AC12").Value = 41
$ws.Range("AE12").Value = 10
$ws.Range("AF12").Value = 17
$ws.Range("AG12").Value = 12
$ws.Range("AH12").Value = 34
$ws.Range("AI12").Value = 26
$ws.Range("G14").Value = 2.1
$ws.Range("H14").Value = 3.4
$ws.Range("I14").Value = 3.5
$ws.Range("J14").Value = 1.06
$ws.Range("K14").Value = 10
$ws.Range("L14").Value = 1.29
$ws.Range("M14").Value = 3.5
$ws.Range("N14").Value = 2
$ws.Range("O14").Value = 1.85
$ws.Range("P14").Value = 1.4
$ws.Range("Q14").Value = 2.75
$ws.Range("R14").Value = 1.8
$ws.Range("S14").Value = 1.95
$ws.Range("T14").Value = 7.5
$ws.Range("U14").Value = 10
$ws.Range("V14").Value = 9
$ws.Range("W14").Value = 19
$ws.Range("X14").Value = 17
$ws.Range("Y14").Value = 26
$ws.Range("Z14").Value = 10
$ws.Range("AA14").Value = 6.5
$ws.Range("AC14").Value = 51
$ws.Range("AD14").Value = 251
$ws.Range("AE14").Value = 10
$ws.Range("AF14").Value = 17
$ws.Range("AG14").Value = 12
$ws.Range("AH14").Value = 41
$ws.Range("AI14").Value = 29
$ws.Range("AJ14").Value = 34
$ws.Range("G15").Value = 1.85
$ws.Range("I15").Value = 3.7
$ws.Range("L15").Value = 1.26
$ws.Range("M15").Value = 3.5
$ws.Range("O15").Value = 1.82
$ws.Range("R15").Value = 1.74
$ws.Range("S15").Value = 1.98
$ws.Range("U15").Value = 7.6
$ws.Range("V15").Value = 7
$ws.Range("W15").Value = 13
$ws.Range("AB15").Value = 11.5
$ws.Range("AF15").Value = 17
$ws.Range("AG15").Value = 10.5
$ws.Range("AI15").Value = 26
$ws.Range("AJ15").Value = 29
$ws.Range("J17").Value = 1.07
$ws.Range("K17").Value = 9
$ws.Range("N17").Value = 2.15
$ws.Range("O17").Value = 1.67
$ws.Range("G18").Value = 1.27
$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 9.5
$ws.Range("J18").Value = 1.02
$ws.Range("K18").Value = 19
$ws.Range("R18").Value = 1.83
$ws.Range("S18").Value = 1.83
$ws.Range("W18").Value = 8.5
$ws.Range("AF18").Value = 51
$ws.Range("G20").Value = 3.6
$ws.Range("H20").Value = 3.2
$ws.Range("I20").Value = 2
$ws.Range("K20").Value = 9.5
$ws.Range("Z20").Value = 9.5
$ws.Range("AJ20").Value = 29
$ws.Range("G22").Value = 1.55
$ws.Range("H22").Value = 3.8
$ws.Range("I22").Value = 5.8
$ws.Range("L22").Value = 1.29
$ws.Range("M22").Value = 3.25
$ws.Range("N22").Value = 1.87
$ws.Range("O22").Value = 1.83
$ws.Range("Q22").Value = 2.67
$ws.Range("U22").Value = 7
$ws.Range("W22").Value = 10.75
$ws.Range("X22").Value = 12.5
$ws.Range("AA22").Value = 7.4
$ws.Range("AB22").Value = 18
$ws.Range("AE22").Value = 14
$ws.Range("AF22").Value = 35
$ws.Range("AG22").Value = 18.5
$ws.Range("AH22").Value = 120
$ws.Range("AI22").Value = 65
$ws.Range("AJ22").Value = 65
$ws.Range("G23").Value = 1.87
$ws.Range("H23").Value = 3.5
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 1.05
$ws.Range("K23").Value = 11
$ws.Range("L23").Value = 1.29
$ws.Range("M23").Value = 3.5
$ws.Range("N23").Value = 2
$ws.Range("O23").Value = 1.85
$ws.Range("P23").Value = 1.4
$ws.Range("Q23").Value = 2.75
$ws.Range("R23").Value = 1.77
$ws.Range("S23").Value = 1.87
$ws.Range("T23").Value = 7
$ws.Range("U23").Value = 9
$ws.Range("V23").Value = 8.5
$ws.Range("W23").Value = 17
$ws.Range("X23").Value = 15
$ws.Range("Y23").Value = 26
$ws.Range("Z23").Value = 10
$ws.Range("AE23").Value = 11
$ws.Range("AF23").Value = 21
$ws.Range("AG23").Value = 13
$ws.Range("AI23").Value = 34
$ws.Range("AJ23").Value = 41
